# Apply updated dSF (column F) values as per the commit:
# "repull data, push all data, mean calculation"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "F3"  = -6
    "F5"  = -6
    "F6"  = -7
    "F7"  = -3
    "F8"  = -5
    "F9"  = -2
    "F10" = -5
    "F11" = -2
    "F12" = -1
    "F13" = -4
    "F14" = 5
    "F15" = 1
    "F16" = 1
    "F17" = 0
    "F19" = -2
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
